$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two JIRA hyperlinks (and their displayed URL text) in A2:A3 were
# removed, leaving the cells blank but still present in the used range.
$target = $ws.Range("A2:A3")
$target.Hyperlinks.Delete()
$target.Clear()

# Selection ends up on A2:A3 (A3 being the active cell).
$ws.Range("A2:A3").Select()
